$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Thomas")

# Add the new log entry in row 28
$ws.Range("A28").Value = "Combining dataframes in python"
$ws.Range("B27").Copy()
$ws.Range("B28").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B28").Value = 44973
$ws.Range("C28").Value = 3

# Extend the SUM formula range to cover the new rows
$ws.Range("E3").Formula = "=SUM(C2:C35)"

$wb.Save()
